$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "26.842.69"
$ws.Range("E2").Value = "  -1.21%  "

Set-TextValue "D3" "1.857.13"
$ws.Range("E3").Value = "  -0.50%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue "D5" "304.31"
$ws.Range("E5").Value = "  -0.97%  "

$ws.Range("E6").Value = "  -0.02%  "

Set-TextValue "D7" "0.5036"
$ws.Range("E7").Value = "  -2.35%  "

Set-TextValue "D8" "0.3647"
$ws.Range("E8").Value = "  -2.78%  "

Set-TextValue "D9" "0.07166"
$ws.Range("E9").Value = "  -0.09%  "

Set-TextValue "D10" "0.8915"
$ws.Range("E10").Value = "  +0.79%  "

Set-TextValue "D11" "20.66"
$ws.Range("E11").Value = "  -0.05%  "

Set-TextValue "D12" "0.07513"
$ws.Range("E12").Value = "  -0.71%  "

Set-TextValue "D13" "1.862.46"
$ws.Range("E13").Value = "  -0.29%  "

Set-TextValue "D14" "92.18"
$ws.Range("E14").Value = "  +3.30%  "

Set-TextValue "D15" "5.227"
$ws.Range("E15").Value = "  -2.05%  "

$ws.Range("E16").Value = "  -0.04%  "

Set-TextValue "D17" "0.000008491"
$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("E18").Value = "  -0.60%  "

$ws.Range("E19").Value = "  +0.00%  "

Set-TextValue "D20" "26.882.26"
$ws.Range("E20").Value = "  -1.23%  "

Set-TextValue "D21" "5.025"
$ws.Range("E21").Value = "  -0.12%  "

Set-TextValue "D22" "2.079.51"
$ws.Range("E22").Value = "  -1.89%  "

$ws.Range("E23").Value = "  -2.59%  "

Set-TextValue "D24" "6.455"
$ws.Range("E24").Value = "  -0.30%  "

Set-TextValue "D25" "146.39"
$ws.Range("E25").Value = "  -3.09%  "

Set-TextValue "D26" "1.793"
$ws.Range("E26").Value = "  -3.08%  "

Set-TextValue "D27" "17.80"
$ws.Range("E27").Value = "  -1.22%  "

Set-TextValue "D28" "2.062"
$ws.Range("E28").Value = "  -3.97%  "

Set-TextValue "D29" "112.76"
$ws.Range("E29").Value = "  -0.13%  "

Set-TextValue "D30" "4.632"
$ws.Range("E30").Value = "  -2.49%  "

Set-TextValue "D31" "4.659"

Set-TextValue "D32" "0.09196"
$ws.Range("E32").Value = "  +2.10%  "

$ws.Range("E33").Value = "  -1.41%  "

Set-TextValue "D34" "2.976"
$ws.Range("E34").Value = "  -4.11%  "

Set-TextValue "D35" "0.7389"
$ws.Range("E35").Value = "  -1.83%  "

Set-TextValue "D36" "1.146"
$ws.Range("E36").Value = "  -2.22%  "

$ws.Range("E37").Value = "  +6.67%  "

Set-TextValue "D38" "2.518"
$ws.Range("E38").Value = "  -0.54%  "

$ws.Range("E39").Value = "  -2.01%  "

Set-TextValue "D40" "1.077"
$ws.Range("E40").Value = "  -0.36%  "

Set-TextValue "D41" "0.5316"
$ws.Range("E41").Value = "  -0.59%  "

Set-TextValue "D42" "119.57"
$ws.Range("E42").Value = "  +3.91%  "

Set-TextValue "D43" "6.476"
$ws.Range("E43").Value = "  -2.57%  "

Set-TextValue "D44" "8.367"
$ws.Range("E44").Value = "  -1.41%  "

Set-TextValue "D45" "0.1463"
$ws.Range("E45").Value = "  -1.40%  "

Set-TextValue "D46" "0.4637"
$ws.Range("E46").Value = "  -0.58%  "

Set-TextValue "D47" "0.9998"
$ws.Range("E47").Value = "  -0.08%  "

Set-TextValue "D48" "9.928"
$ws.Range("E48").Value = "  -1.89%  "

Set-TextValue "D49" "1.560"
$ws.Range("E49").Value = "  -0.76%  "

Set-TextValue "D50" "36.89"
$ws.Range("E50").Value = "  +1.33%  "

Set-TextValue "D51" "62.78"
$ws.Range("E51").Value = "  -3.34%  "
